# daily auto push: 2026-01-18 06:46 UTC
# A new observation (2026/01/18, 日, 13, 194) was recorded and inserted
# into the time series at its sorted position (row 674), pushing every
# following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 674..715 down to 675..716, opening up a blank row 674.
$ws.Rows.Item(674).Insert()

# Column A holds dates stored as literal text (e.g. "2026/12/29"), not
# real Excel date serials. Force the cell to Text format *before* writing
# the value so the slash-separated string isn't auto-coerced into a date.
$ws.Cells.Item(674, 1).NumberFormat = "@"
$ws.Cells.Item(674, 1).Value = "2026/01/18"
$ws.Cells.Item(674, 2).Value = "日"
$ws.Cells.Item(674, 3).Value = 13
$ws.Cells.Item(674, 4).Value = 194
